# Adding test cases for Notifications
# - "Test Cases" sheet (sheet1): row 20 (TestCase_F19) Results flips PASS -> SKIP,
#   and two brand-new rows (21, 22) are appended for TestCase_F20 / TestCase_F21.
# - The sheet view scrolls right a bit and the selection moves to D15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing row 20: Results PASS -> SKIP -------------------------------
$ws.Range("E20").Value = "SKIP"

# --- new row 21: TestCase_F20 --------------------------------------------
# Copy the formatting of the row above (row 19) onto row 21 first so the new
# row picks up the same borders/fill as its neighbours, then fill in values.
$ws.Range("A19:E19").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)

$ws.Range("A21").Value = "TestCase_F20"
$ws.Range("B21").Value = "OPQA-1011"
$ws.Range("C21").Value = "Verify that follower of the post is able to start conversation from home page when some one commented on the post he is following."
$ws.Range("D21").Value = "Y"
$ws.Range("E21").Value = "SKIP"

# --- new row 22: TestCase_F21 --------------------------------------------
$ws.Range("A19:E19").Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)

$ws.Range("A22").Value = "TestCase_F21"
$ws.Range("B22").Value = "OPQA-1010"
$ws.Range("C22").Value = "Verify that author of the post is able to start conversation from home page when some one commented on his post."
$ws.Range("D22").Value = "Y"
$ws.Range("E22").Value = "PASS"

# --- view state: scroll right one column, move the active selection ------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("D15").Select()
